# resumefit-tracker.xlsx -- v0.0.3: limit characters, add download endpoint,
# versioning, cookie tracking, release notes.
#
# Applies the row/value edits on the "features" sheet:
#  - F11 gets its missing 0.0.3 release date
#  - E12/F12 gain the 0.0.3 version + date (previously blank)
#  - B17 typo fix: "asekd" -> "asked"
#  - E19/F19 gain the 0.0.3 version + date (previously blank)
#  - three brand-new feature rows (21-23): domain redirect, release notes,
#    tab icon update
#  - two new trailing rows (25-26) continuing the "sl no" counter
#  - page setup (paper size / orientation) + selection get refreshed

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("features")
$ws.Activate()

# --- row 11: backfill the release date that was missing ---
$ws.Range("F11").Value = 45577

# --- row 12: backfill version + date (copy date formatting from a sibling
#     date cell so the new cell picks up the existing numFmt style instead
#     of manufacturing a new one) ---
$ws.Range("E12").Value = "0.0.3"
$ws.Range("F13").Copy()
$ws.Range("F12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F12").Value = 45577

# --- row 17: fix the "asekd" -> "asked" typo ---
$ws.Range("B17").Value = "top 10 questions you may be asked with answers"

# --- row 19: backfill version + date ---
$ws.Range("E19").Value = "0.0.3"
$ws.Range("F13").Copy()
$ws.Range("F19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F19").Value = 45577

$excel.CutCopyMode = $false

# --- new rows: 21, 22, 23 ---
$ws.Range("B21").Value = "resolve domain redirect"
$ws.Range("C21").Value = "shiv"
$ws.Range("D21").Value = 1

$ws.Range("B22").Value = "add release notes to keep track "
$ws.Range("C22").Value = "vivek"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = "0.0.3"

$ws.Range("B23").Value = "tab icon update from vue"
$ws.Range("C23").Value = "vivek"
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = "0.0.3"

$ws.Range("F13").Copy()
$ws.Range("F22:F23").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F22").Value = 45577
$ws.Range("F23").Value = 45577

# --- two more trailing "sl no" rows ---
$ws.Range("A25").Value = 24
$ws.Range("A26").Value = 25

# --- page setup: paper size + portrait orientation ---
$ws.PageSetup.PaperSize = 9        # xlPaperA4
$ws.PageSetup.Orientation = 1      # xlPortrait

# --- refresh selection to match the saved view ---
$ws.Range("G10").Select()

Write-Output "resumefit-tracker v0.0.3 edits applied"
